$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-14 Sunday" "2025-09-15 Monday"

Replace-Text "636×8=" "103×8="
Replace-Text "107×8=" "854×9="
Replace-Text "430×6=" "583×7="
Replace-Text "370×2=" "404×2="
Replace-Text "296×4=" "103×4="
Replace-Text "500×3=" "754×8="
Replace-Text "368×7=" "691×5="
Replace-Text "484×4=" "543×6="
Replace-Text "385×2=" "391×2="
Replace-Text "374×9=" "336×9="
Replace-Text "475×9=" "645×9="
Replace-Text "476×5=" "565×6="
Replace-Text "193×6=" "799×2="
Replace-Text "564×8=" "146×2="
Replace-Text "248×3=" "172×3="
Replace-Text "738×7=" "288×2="
Replace-Text "315×2=" "123×8="
Replace-Text "114×6=" "155×9="
Replace-Text "147×6=" "181×7="
Replace-Text "393×3=" "172×5="
Replace-Text "822×5=" "921×5="
Replace-Text "807×3=" "232×6="
Replace-Text "691×9=" "987×7="
Replace-Text "611×3=" "874×5="
Replace-Text "222×5=" "493×2="
